$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F column values (想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 199
$ws1.Range("F4").Value = 429
$ws1.Range("F5").Value = 471
$ws1.Range("F6").Value = 277
$ws1.Range("F7").Value = 2507
$ws1.Range("F9").Value = 6681

# Sheet "全部类型" (All types) - same updates, mirrored with different row for last item
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 199
$ws4.Range("F4").Value = 429
$ws4.Range("F5").Value = 471
$ws4.Range("F6").Value = 277
$ws4.Range("F9").Value = 2507
$ws4.Range("F11").Value = 6681
